# Update the regional swings (polling figures) on the "Calc" sheet.
# A new poll's data was inserted at row 6 ("Latest Morgan"), which pushed the
# previously-latest figures down into rows 7 ("Second Morgan") and 8
# ("Third Morgan").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Row 6 - Latest Morgan (new poll data)
$ws.Range("B6").Value = 56.5
$ws.Range("C6").Value = 55.5
$ws.Range("D6").Value = 58.5
$ws.Range("E6").Value = 54.5
$ws.Range("F6").Value = 50.5
$ws.Range("G6").Value = 64.5

# Row 7 - Second Morgan (previous "Latest Morgan" figures)
$ws.Range("B7").Value = 55.5
$ws.Range("C7").Value = 55.5
$ws.Range("D7").Value = 58
$ws.Range("E7").Value = 51.5
$ws.Range("F7").Value = 53.5
$ws.Range("G7").Value = 55.5

# Row 8 - Third Morgan (previous "Second Morgan" figures)
$ws.Range("B8").Value = 53.5
$ws.Range("C8").Value = 53.5
$ws.Range("D8").Value = 55
$ws.Range("E8").Value = 47
$ws.Range("F8").Value = 53.5
$ws.Range("G8").Value = 57.5

$excel.Calculate()

# Update the active selection to reflect where the user ended up working.
$ws.Range("H30").Select()
